$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-17 (16 rows) x columns A-T (20 cols) of the LR-pairs result table,
# recomputed per the new natmi clustering (adds the sCs->* rows and refreshes
# all specificity/weight figures).
$rows = @(
  @('ECs', 'Apod', 'Lepr', 'ECs', 3, 1, 1.138252333333333, 3.414757, 0.004456632015781552, 0.004456632015781552, 3, 1, 2.473323333333334, 7.41997, 0.5382827202858395, 0.5382827202858395, 2.815266055254444, 25.33739449729, 0.002398928004767858, 0.002398928004767858),
  @('ECs', 'Apod', 'Lepr', 'FAPs', 3, 1, 1.138252333333333, 3.414757, 0.004456632015781552, 0.004456632015781552, 3, 1, 1.551135333333334, 4.653406, 0.3375819633063809, 0.3375819633063809, 1.765583412482445, 15.890250712342, 0.00150447858562161, 0.00150447858562161),
  @('ECs', 'Apod', 'Lepr', 'M2', 3, 1, 1.138252333333333, 3.414757, 0.004456632015781552, 0.004456632015781552, 1, 0.3333333333333333, 0.01982566666666667, 0.059477, 0.004314766953834163, 0.004314766953834163, 0.02256661134322222, 0.203099502089, [double]"1.922932854709357E-05", [double]"1.922932854709357E-05"),
  @('ECs', 'Apod', 'Lepr', 'sCs', 3, 1, 1.138252333333333, 3.414757, 0.004456632015781552, 0.004456632015781552, 3, 1, 0.5505563333333333, 1.651669, 0.1198205494539455, 0.1198205494539455, 0.626672031048111, 5.640048279433, 0.0005339960968449903, 0.0005339960968449904),
  @('FAPs', 'Apod', 'Lepr', 'ECs', 3, 1, 245.1331126666666, 735.3993379999999, 0.9597767085960608, 0.9597767085960607, 3, 1, 2.473323333333334, 7.41997, 0.5382827202858395, 0.5382827202858395, 606.2934473310955, 5456.641025979859, 0.516631217570077, 0.516631217570077),
  @('FAPs', 'Apod', 'Lepr', 'FAPs', 3, 1, 245.1331126666666, 735.3993379999999, 0.9597767085960608, 0.9597767085960607, 3, 1, 1.551135333333334, 4.653406, 0.3375819633063809, 0.3375819633063809, 380.2346324272476, 3422.111691845228, 0.3240033056235944, 0.3240033056235944),
  @('FAPs', 'Apod', 'Lepr', 'M2', 3, 1, 245.1331126666666, 735.3993379999999, 0.9597767085960608, 0.9597767085960607, 1, 0.3333333333333333, 0.01982566666666667, 0.059477, 0.004314766953834163, 0.004314766953834163, 4.859927380691778, 43.739346426226, 0.004141212825310004, 0.004141212825310004),
  @('FAPs', 'Apod', 'Lepr', 'sCs', 3, 1, 245.1331126666666, 735.3993379999999, 0.9597767085960608, 0.9597767085960607, 3, 1, 0.5505563333333333, 1.651669, 0.1198205494539455, 0.1198205494539455, 134.9595876883469, 1214.636289195122, 0.1150009725770793, 0.1150009725770794),
  @('M2', 'Apod', 'Lepr', 'ECs', 3, 1, 1.996800333333334, 5.990401, 0.007818129630884373, 0.007818129630884373, 3, 1, 2.473323333333334, 7.41997, 0.5382827202858395, 0.5382827202858395, 4.938732856441112, 44.44859570797, 0.004208364085259766, 0.004208364085259766),
  @('M2', 'Apod', 'Lepr', 'FAPs', 3, 1, 1.996800333333334, 5.990401, 0.007818129630884373, 0.007818129630884373, 3, 1, 1.551135333333334, 4.653406, 0.3375819633063809, 0.3375819633063809, 3.097307550645112, 27.875767955806, 0.002639259550177738, 0.002639259550177738),
  @('M2', 'Apod', 'Lepr', 'M2', 3, 1, 1.996800333333334, 5.990401, 0.007818129630884373, 0.007818129630884373, 1, 0.3333333333333333, 0.01982566666666667, 0.059477, 0.004314766953834163, 0.004314766953834163, 0.03958789780855556, 0.356291080277, [double]"3.373340737213158E-05", [double]"3.373340737213158E-05"),
  @('M2', 'Apod', 'Lepr', 'sCs', 3, 1, 1.996800333333334, 5.990401, 0.007818129630884373, 0.007818129630884373, 3, 1, 0.5505563333333333, 1.651669, 0.1198205494539455, 0.1198205494539455, 1.099351069918778, 9.894159629269001, 0.0009367725880747377, 0.0009367725880747379),
  @('sCs', 'Apod', 'Lepr', 'ECs', 3, 1, 7.138233333333333, 21.4147, 0.02794852975727327, 0.02794852975727327, 3, 1, 2.473323333333334, 7.41997, 0.5382827202858395, 0.5382827202858395, 17.65515906211111, 158.896431559, 0.01504421062573479, 0.01504421062573479),
  @('sCs', 'Apod', 'Lepr', 'FAPs', 3, 1, 7.138233333333333, 21.4147, 0.02794852975727327, 0.02794852975727327, 3, 1, 1.551135333333334, 4.653406, 0.3375819633063809, 0.3375819633063809, 11.07236594091111, 99.65129346820001, 0.009434919546987119, 0.009434919546987119),
  @('sCs', 'Apod', 'Lepr', 'M2', 3, 1, 7.138233333333333, 21.4147, 0.02794852975727327, 0.02794852975727327, 1, 0.3333333333333333, 0.01982566666666667, 0.059477, 0.004314766953834163, 0.004314766953834163, 0.1415202346555555, 1.2736821119, 0.0001205913926049334, 0.0001205913926049334),
  @('sCs', 'Apod', 'Lepr', 'sCs', 3, 1, 7.138233333333333, 21.4147, 0.02794852975727327, 0.02794852975727327, 3, 1, 0.5505563333333333, 1.651669, 0.1198205494539455, 0.1198205494539455, 3.929999570477777, 35.3699961343, 0.003348808191946429, 0.00334880819194643)
)

$data = New-Object 'object[,]' 16,20
for ($i = 0; $i -lt 16; $i++) {
  for ($j = 0; $j -lt 20; $j++) {
    $data[$i,$j] = $rows[$i][$j]
  }
}

$ws.Range("A2:T17").Value = $data

Write-Output "applied"
